$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 33.47808166666666
$ws.Range("H2").Value = 100.434245
$ws.Range("I2").Value = 0.4880542983452505
$ws.Range("J2").Value = 0.4880542983452505
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 55.908252
$ws.Range("N2").Value = 167.724756
$ws.Range("O2").Value = 0.6412441619121594
$ws.Range("P2").Value = 0.6412441619121594
$ws.Range("Q2").Value = 1871.70102629658
$ws.Range("R2").Value = 16845.30923666922
$ws.Range("S2").Value = 0.3129619695100271
$ws.Range("T2").Value = 0.3129619695100271

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 33.47808166666666
$ws.Range("H3").Value = 100.434245
$ws.Range("I3").Value = 0.4880542983452505
$ws.Range("J3").Value = 0.4880542983452505
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.370676
$ws.Range("N3").Value = 16.112028
$ws.Range("O3").Value = 0.06159939735768789
$ws.Range("P3").Value = 0.06159939735768789
$ws.Range("Q3").Value = 179.7999297332067
$ws.Range("R3").Value = 1618.19936759886
$ws.Range("S3").Value = 0.03006385065589664
$ws.Range("T3").Value = 0.03006385065589664

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 33.47808166666666
$ws.Range("H4").Value = 100.434245
$ws.Range("I4").Value = 0.4880542983452505
$ws.Range("J4").Value = 0.4880542983452505
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.90822366666667
$ws.Range("N4").Value = 77.724671
$ws.Range("O4").Value = 0.2971564407301527
$ws.Range("P4").Value = 0.2971564407301527
$ws.Range("Q4").Value = 867.3576277509327
$ws.Range("R4").Value = 7806.218649758394
$ws.Range("S4").Value = 0.1450284781793267
$ws.Range("T4").Value = 0.1450284781793267

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 23.24776266666667
$ws.Range("H5").Value = 69.743288
$ws.Range("I5").Value = 0.3389134003957588
$ws.Range("J5").Value = 0.3389134003957588
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 55.908252
$ws.Range("N5").Value = 167.724756
$ws.Range("O5").Value = 0.6412441619121594
$ws.Range("P5").Value = 0.6412441619121594
$ws.Range("Q5").Value = 1299.741773604192
$ws.Range("R5").Value = 11697.67596243773
$ws.Range("S5").Value = 0.2173262393975785
$ws.Range("T5").Value = 0.2173262393975784

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.24776266666667
$ws.Range("H6").Value = 69.743288
$ws.Range("I6").Value = 0.3389134003957588
$ws.Range("J6").Value = 0.3389134003957588
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.370676
$ws.Range("N6").Value = 16.112028
$ws.Range("O6").Value = 0.06159939735768789
$ws.Range("P6").Value = 0.06159939735768789
$ws.Range("Q6").Value = 124.8562010075627
$ws.Range("R6").Value = 1123.705809068064
$ws.Range("S6").Value = 0.02087686122082352
$ws.Range("T6").Value = 0.02087686122082352

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 23.24776266666667
$ws.Range("H7").Value = 69.743288
$ws.Range("I7").Value = 0.3389134003957588
$ws.Range("J7").Value = 0.3389134003957588
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 25.90822366666667
$ws.Range("N7").Value = 77.724671
$ws.Range("O7").Value = 0.2971564407301527
$ws.Range("P7").Value = 0.2971564407301527
$ws.Range("Q7").Value = 602.3082349175833
$ws.Range("R7").Value = 5420.774114258249
$ws.Range("S7").Value = 0.1007102997773568
$ws.Range("T7").Value = 0.1007102997773568

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.86914966666667
$ws.Range("H8").Value = 35.607449
$ws.Range("I8").Value = 0.1730323012589908
$ws.Range("J8").Value = 0.1730323012589908
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 55.908252
$ws.Range("N8").Value = 167.724756
$ws.Range("O8").Value = 0.6412441619121594
$ws.Range("P8").Value = 0.6412441619121594
$ws.Range("Q8").Value = 663.583410589716
$ws.Range("R8").Value = 5972.250695307443
$ws.Range("S8").Value = 0.1109559530045539
$ws.Range("T8").Value = 0.1109559530045538

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.86914966666667
$ws.Range("H9").Value = 35.607449
$ws.Range("I9").Value = 0.1730323012589908
$ws.Range("J9").Value = 0.1730323012589908
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.370676
$ws.Range("N9").Value = 16.112028
$ws.Range("O9").Value = 0.06159939735768789
$ws.Range("P9").Value = 0.06159939735768789
$ws.Range("Q9").Value = 63.74535725517466
$ws.Range("R9").Value = 573.708215296572
$ws.Range("S9").Value = 0.01065868548096773
$ws.Range("T9").Value = 0.01065868548096773

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.86914966666667
$ws.Range("H10").Value = 35.607449
$ws.Range("I10").Value = 0.1730323012589908
$ws.Range("J10").Value = 0.1730323012589908
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 25.90822366666667
$ws.Range("N10").Value = 77.724671
$ws.Range("O10").Value = 0.2971564407301527
$ws.Range("P10").Value = 0.2971564407301527
$ws.Range("Q10").Value = 307.5085842971421
$ws.Range("R10").Value = 2767.577258674279
$ws.Range("S10").Value = 0.05141766277346923
$ws.Range("T10").Value = 0.05141766277346922
